# Model_Optimization_Notes.xlsx - add linear regression / random forest
# regressor trial rows (Attempt IDs 6-11) plus a new "n_estimators" column,
# matching the author's updated model trial data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear out the existing data block (values AND formatting) so the old
#     column layout (which had no "n_estimators" column) doesn't leave stray
#     cells/styles behind once the columns shift right by one. -------------
$ws.Range("A1:K20").Clear()

# --- Header row -------------------------------------------------------------
$ws.Cells.Item(1,1).Value  = "Attempt ID"
$ws.Cells.Item(1,2).Value  = "Model Type"
$ws.Cells.Item(1,3).Value  = "Features"
$ws.Cells.Item(1,4).Value  = "Data Entries"
$ws.Cells.Item(1,5).Value  = "Total Params"
$ws.Cells.Item(1,6).Value  = "n_estimators"
$ws.Cells.Item(1,7).Value  = "Loss Function"
$ws.Cells.Item(1,8).Value  = "Loss"
$ws.Cells.Item(1,9).Value  = "Target Stdev"
$ws.Cells.Item(1,10).Value = "RMSE"
$ws.Cells.Item(1,11).Value = "r^2 value"

# --- Row 2 : Attempt 0 (neural network, aborted) ---------------------------
$ws.Cells.Item(2,1).Value  = 0
$ws.Cells.Item(2,2).Value  = "neural network"
$ws.Cells.Item(2,3).Value  = 89
$ws.Cells.Item(2,4).Value  = 2537
$ws.Cells.Item(2,5).Value  = 6025
$ws.Cells.Item(2,7).Value  = "BCE"
$ws.Cells.Item(2,8).Value  = "aborted"
$ws.Cells.Item(2,9).Value  = "n/a"
$ws.Cells.Item(2,10).Value = "n/a"
$ws.Cells.Item(2,11).Value = "n/a"

# --- Row 3 : Attempt 1 (neural network, MAE) --------------------------------
$ws.Cells.Item(3,1).Value  = 1
$ws.Cells.Item(3,2).Value  = "neural network"
$ws.Cells.Item(3,3).Value  = 89
$ws.Cells.Item(3,4).Value  = 2537
$ws.Cells.Item(3,5).Value  = 6025
$ws.Cells.Item(3,7).Value  = "MAE"
$ws.Cells.Item(3,8).Value  = 2.632
$ws.Cells.Item(3,9).Value  = "n/a"
$ws.Cells.Item(3,10).Value = "n/a"
$ws.Cells.Item(3,11).Value = "n/a"

# --- Row 4 : Attempt 2 (neural network, MSE) --------------------------------
$ws.Cells.Item(4,1).Value  = 2
$ws.Cells.Item(4,2).Value  = "neural network"
$ws.Cells.Item(4,3).Value  = 89
$ws.Cells.Item(4,4).Value  = 2537
$ws.Cells.Item(4,5).Value  = 8113
$ws.Cells.Item(4,7).Value  = "MSE"
$ws.Cells.Item(4,8).Value  = 15.17
$ws.Cells.Item(4,9).Value  = "n/a"
$ws.Cells.Item(4,10).Value = "n/a"
$ws.Cells.Item(4,11).Value = 0.6875

# --- Row 5 : Attempt 3 (neural network, MSE) --------------------------------
$ws.Cells.Item(5,1).Value  = 3
$ws.Cells.Item(5,2).Value  = "neural network"
$ws.Cells.Item(5,3).Value  = 89
$ws.Cells.Item(5,4).Value  = 2537
$ws.Cells.Item(5,5).Value  = 8417
$ws.Cells.Item(5,7).Value  = "MSE"
$ws.Cells.Item(5,8).Value  = "n/a"
$ws.Cells.Item(5,9).Value  = "n/a"
$ws.Cells.Item(5,10).Value = "n/a"
$ws.Cells.Item(5,11).Value = 0.7457
# This cell historically carries a distinct font colour (style index 1).
$ws.Cells.Item(5,11).Font.Color = 1907741

# --- Row 6 : Attempt 4 (neural network, MSE) --------------------------------
$ws.Cells.Item(6,1).Value  = 4
$ws.Cells.Item(6,2).Value  = "neural network"
$ws.Cells.Item(6,3).Value  = 87
$ws.Cells.Item(6,4).Value  = 6289
$ws.Cells.Item(6,5).Value  = 8417
$ws.Cells.Item(6,7).Value  = "MSE"
$ws.Cells.Item(6,8).Value  = 9.0807
$ws.Cells.Item(6,9).Value  = "n/a"
$ws.Cells.Item(6,10).Value = "n/a"
$ws.Cells.Item(6,11).Value = 0.8324

# --- Row 7 : Attempt 5 (neural network, MSE) --------------------------------
$ws.Cells.Item(7,1).Value  = 5
$ws.Cells.Item(7,2).Value  = "neural network"
$ws.Cells.Item(7,3).Value  = 87
$ws.Cells.Item(7,4).Value  = 6289
$ws.Cells.Item(7,5).Value  = 8417
$ws.Cells.Item(7,7).Value  = "MSE"
$ws.Cells.Item(7,8).Value  = 8.7319
$ws.Cells.Item(7,9).Value  = 7.1963
$ws.Cells.Item(7,10).Value = 2.955
$ws.Cells.Item(7,11).Value = 0.8312

# --- Row 8 : Attempt 6 (linear regression) ----------------------------------
$ws.Cells.Item(8,1).Value  = 6
$ws.Cells.Item(8,2).Value  = "linear regression"
$ws.Cells.Item(8,3).Value  = 87
$ws.Cells.Item(8,4).Value  = 6289
$ws.Cells.Item(8,5).Value  = "n/a"
$ws.Cells.Item(8,7).Value  = "MSE"
$ws.Cells.Item(8,8).Value  = 14.943
$ws.Cells.Item(8,9).Value  = 7.6946
$ws.Cells.Item(8,10).Value = 3.8656
$ws.Cells.Item(8,11).Value = 0.743

# --- Row 9 : Attempt 7 (random forest regressor, n_estimators=25) ----------
$ws.Cells.Item(9,1).Value  = 7
$ws.Cells.Item(9,2).Value  = "random forest regressor"
$ws.Cells.Item(9,3).Value  = 87
$ws.Cells.Item(9,4).Value  = 6289
$ws.Cells.Item(9,5).Value  = "n/a"
$ws.Cells.Item(9,6).Value  = 25
$ws.Cells.Item(9,7).Value  = "MSE"
$ws.Cells.Item(9,8).Value  = 13.4441
$ws.Cells.Item(9,9).Value  = 7.2005
$ws.Cells.Item(9,10).Value = 3.6666
$ws.Cells.Item(9,11).Value = 0.7404

# --- Row 10 : Attempt 8 (random forest regressor, n_estimators=50) ---------
$ws.Cells.Item(10,1).Value  = 8
$ws.Cells.Item(10,2).Value  = "random forest regressor"
$ws.Cells.Item(10,3).Value  = 87
$ws.Cells.Item(10,4).Value  = 6289
$ws.Cells.Item(10,5).Value  = "n/a"
$ws.Cells.Item(10,6).Value  = 50
$ws.Cells.Item(10,7).Value  = "MSE"
$ws.Cells.Item(10,8).Value  = 12.8613
$ws.Cells.Item(10,9).Value  = 7.2005
$ws.Cells.Item(10,10).Value = 3.5862
$ws.Cells.Item(10,11).Value = 0.7517

# --- Row 11 : Attempt 9 (random forest regressor, n_estimators=100) --------
$ws.Cells.Item(11,1).Value  = 9
$ws.Cells.Item(11,2).Value  = "random forest regressor"
$ws.Cells.Item(11,3).Value  = 87
$ws.Cells.Item(11,4).Value  = 6289
$ws.Cells.Item(11,5).Value  = "n/a"
$ws.Cells.Item(11,6).Value  = 100
$ws.Cells.Item(11,7).Value  = "MSE"
$ws.Cells.Item(11,8).Value  = 12.7365
$ws.Cells.Item(11,9).Value  = 7.2005
$ws.Cells.Item(11,10).Value = 3.5688
$ws.Cells.Item(11,11).Value = 0.7541

# --- Row 12 : Attempt 10 (random forest regressor, n_estimators=250) -------
$ws.Cells.Item(12,1).Value  = 10
$ws.Cells.Item(12,2).Value  = "random forest regressor"
$ws.Cells.Item(12,3).Value  = 87
$ws.Cells.Item(12,4).Value  = 6289
$ws.Cells.Item(12,5).Value  = "n/a"
$ws.Cells.Item(12,6).Value  = 250
$ws.Cells.Item(12,7).Value  = "MSE"
$ws.Cells.Item(12,8).Value  = 12.6924
$ws.Cells.Item(12,9).Value  = 7.2005
$ws.Cells.Item(12,10).Value = 3.5626
$ws.Cells.Item(12,11).Value = 0.755

# --- Row 13 : Attempt 11 (random forest regressor, n_estimators=1000) ------
$ws.Cells.Item(13,1).Value  = 11
$ws.Cells.Item(13,2).Value  = "random forest regressor"
$ws.Cells.Item(13,3).Value  = 87
$ws.Cells.Item(13,4).Value  = 6289
$ws.Cells.Item(13,5).Value  = "n/a"
$ws.Cells.Item(13,6).Value  = 1000
$ws.Cells.Item(13,7).Value  = "MSE"
$ws.Cells.Item(13,8).Value  = 12.6256
$ws.Cells.Item(13,9).Value  = 7.2005
$ws.Cells.Item(13,10).Value = 3.553
$ws.Cells.Item(13,11).Value = 0.7562

# --- Active selection, matching the author's last cursor position ----------
$ws.Range("C9").Select()
